$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.433.89"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.549.25"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.53"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.480"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.97"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "1.770.73"
$ws.Range("D13").Value = "1.554.11"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "28.398.44"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.93"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.26"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Value = "0.0₃0672"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.91"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.64"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.73"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.102"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0466"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("E31").Value = "  -4.52%  "
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").Value = "1.383.49"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("E36").Value = "  -3.17%  "
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.56"
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("E40").Value = "  +1.62%  "
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0457"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.76"
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("D47").Value = "1.683.00"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.873"
$ws.Range("E48").Value = "  -9.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.60"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.91"
$ws.Range("E50").Value = "  +7.73%  "
$ws.Range("E51").Value = "  -2.03%  "
